$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "63.756.46"
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -0.33%  "
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.617.77"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -0.12%  "
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.03%  "
$c.ClearFormats()
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "595.19"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -0.51%  "
$c.ClearFormats()
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "150.86"
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +0.59%  "
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +0.03%  "
$c.ClearFormats()
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.585"
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -0.19%  "
$c.ClearFormats()
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.113"
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +3.65%  "
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +3.38%  "
$c.ClearFormats()
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.393"
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +2.67%  "
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  +0.63%  "
$c.ClearFormats()
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "27.97"
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +1.53%  "
$c.ClearFormats()
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.087.56"
$c.ClearFormats()
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "63.569.21"
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -0.40%  "
$c.ClearFormats()
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000163"
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +9.96%  "
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.602.99"
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -0.71%  "
$c.ClearFormats()
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "12.24"
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +0.55%  "
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +3.19%  "
$c.ClearFormats()
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "347.49"
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  -0.66%  "
$c.ClearFormats()
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.98"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +1.28%  "
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +0.22%  "
$c.ClearFormats()
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "67.29"
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +1.55%  "
$c.ClearFormats()
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.69"
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -2.28%  "
$c.ClearFormats()
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.27"
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +0.70%  "
$c.ClearFormats()
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.69"
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +0.11%  "
$c.ClearFormats()
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.50"
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +4.67%  "
$c.ClearFormats()
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "552.49"
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +1.71%  "
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -0.65%  "
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -0.11%  "
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +0.55%  "
$c.ClearFormats()
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0₃0873"
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +2.87%  "
$c.ClearFormats()
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.79"
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +1.98%  "
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +6.16%  "
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +0.66%  "
$c.ClearFormats()
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "164.18"
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  -2.39%  "
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +2.51%  "
$c.ClearFormats()
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.98"
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  +1.17%  "
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -0.03%  "
$c.ClearFormats()
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "19.67"
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +1.28%  "
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -0.05%  "
$c.ClearFormats()
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "166.85"
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -0.82%  "
$c.ClearFormats()
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "4.09"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  +4.43%  "
$c.ClearFormats()
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "23.56"
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +10.00%  "
$c.ClearFormats()
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0584"
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -0.47%  "
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +7.99%  "
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +1.52%  "
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +2.76%  "
$c.ClearFormats()
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0964"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -0.01%  "
$c.ClearFormats()
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "19.37"
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +1.33%  "
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0₆0230"
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +17.59%  "
$c.ClearFormats()
